$wb = $excel.ActiveWorkbook

# --- Moorings sheet: fill in the missing "Recover Date" calibration event ---
$moorings = $wb.Worksheets.Item("Moorings")
$moorings.Range("G2").Value = "12/20/2015"

# --- Asset_Cal_Info sheet: assign the missing OOI bar code for the ENG sensor ---
$calInfo = $wb.Worksheets.Item("Asset_Cal_Info")
$calInfo.Range("E11").Style = "Normal"
$calInfo.Range("E11").Value = "OL000135"

# --- restore the final selection / active sheet state ---
$calInfo.Range("F20").Select() | Out-Null
$moorings.Activate() | Out-Null
$moorings.Range("G11").Select() | Out-Null
